# Update Visibility_Results sheet (sheet 1): columns A (Satellite_ID), B (IoT_ID), C (Duration (s))
# for rows 2-37.
$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Visibility_Results")

$sheet1Data = @(
    @(2, 0, 0, 518),
    @(3, 0, 0, 196),
    @(4, 0, 2, 579),
    @(5, 0, 2, 136),
    @(6, 0, 3, 458),
    @(7, 0, 3, 256),
    @(8, 0, 5, 459),
    @(9, 0, 5, 256),
    @(10, 0, 6, 457),
    @(11, 0, 6, 260),
    @(12, 0, 9, 578),
    @(13, 0, 9, 135),
    @(14, 0, 8, 11),
    @(15, 0, 8, 323),
    @(16, 0, 8, 265),
    @(17, 0, 1, 395),
    @(18, 0, 1, 260),
    @(19, 0, 4, 449),
    @(20, 0, 4, 259),
    @(21, 0, 7, 926),
    @(22, 1, 7, 916),
    @(23, 1, 9, 131),
    @(24, 1, 9, 520),
    @(25, 1, 2, 77),
    @(26, 1, 2, 522),
    @(27, 1, 0, 15),
    @(28, 1, 0, 584),
    @(29, 1, 3, 4),
    @(30, 1, 3, 643),
    @(31, 1, 5, 4),
    @(32, 1, 5, 644),
    @(33, 1, 8, 264),
    @(34, 1, 8, 265),
    @(35, 1, 1, 699),
    @(36, 1, 4, 697),
    @(37, 1, 6, 646)
)

foreach ($row in $sheet1Data) {
    $r = $row[0]
    $ws1.Cells.Item($r, 1).Value = $row[1]
    $ws1.Cells.Item($r, 2).Value = $row[2]
    $ws1.Cells.Item($r, 3).Value = $row[3]
}

# Update Link_Budget_Results sheet (sheet 2): columns A (Satellite_ID), B (IoT_ID),
# F (Uplink_Bitrate_bps), G (Downlink_Bitrate_bps) for rows 2-21.
$ws2 = $wb.Worksheets.Item("Link_Budget_Results")

$sheet2Data = @(
    @(2, 0, 0, 0.1930501930501931, 19305.0193050193),
    @(3, 0, 2, 0.1727115716753022, 17271.15716753023),
    @(4, 0, 3, 0.2183406113537118, 21834.06113537118),
    @(5, 0, 5, 0.2178649237472767, 21786.49237472767),
    @(6, 0, 6, 0.2188183807439825, 21881.83807439825),
    @(7, 0, 9, 0.1730103806228374, 17301.03806228374),
    @(8, 0, 8, 9.090909090909092, 909090.9090909091),
    @(9, 0, 1, 0.2531645569620253, 25316.45569620253),
    @(10, 0, 4, 0.22271714922049, 22271.714922049),
    @(11, 0, 7, 0.1079913606911447, 10799.13606911447),
    @(12, 1, 7, 0.1091703056768559, 10917.03056768559),
    @(13, 1, 9, 0.7633587786259542, 76335.87786259541),
    @(14, 1, 2, 1.298701298701299, 129870.1298701299),
    @(15, 1, 0, 6.666666666666667, 666666.6666666666),
    @(16, 1, 3, 25, 2500000),
    @(17, 1, 5, 25, 2500000),
    @(18, 1, 8, 0.3787878787878788, 37878.78787878788),
    @(19, 1, 1, 0.1430615164520744, 14306.15164520744),
    @(20, 1, 4, 0.1434720229555237, 14347.20229555237),
    @(21, 1, 6, 0.1547987616099071, 15479.87616099071)
)

foreach ($row in $sheet2Data) {
    $r = $row[0]
    $ws2.Cells.Item($r, 1).Value = $row[1]
    $ws2.Cells.Item($r, 2).Value = $row[2]
    $ws2.Cells.Item($r, 6).Value = $row[3]
    $ws2.Cells.Item($r, 7).Value = $row[4]
}
